# Update the NATMI TPM-derived values on the active sheet (Lama1-Itga7 LR-pair table).
# Only the "Receptor average/total expression value" for the ECs target cluster changed
# (rows 2 and 5, column M/N), and all downstream derived-specificity / edge-weight
# columns (O, P, Q, R, S, T) that depend on those values were recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster FAPs, Target cluster ECs)
$ws.Range("M2").Value = 2.750415333333333
$ws.Range("N2").Value = 8.251245999999998
$ws.Range("O2").Value = 0.04811444325525444
$ws.Range("P2").Value = 0.04811444325525444
$ws.Range("Q2").Value = 0.05431245158733333
$ws.Range("R2").Value = 0.4888120642859999
$ws.Range("S2").Value = 0.03436223909444881
$ws.Range("T2").Value = 0.03436223909444881

# Row 3 (Sending cluster FAPs, Target cluster FAPs)
$ws.Range("O3").Value = 0.01415294505639593
$ws.Range("P3").Value = 0.01415294505639593
$ws.Range("S3").Value = 0.01010771088711213
$ws.Range("T3").Value = 0.01010771088711213

# Row 4 (Sending cluster FAPs, Target cluster MuSCs)
$ws.Range("O4").Value = 0.9377326116883496
$ws.Range("P4").Value = 0.9377326116883496
$ws.Range("R4").Value = 9.526765408797001
$ws.Range("S4").Value = 0.6697072652083125
$ws.Range("T4").Value = 0.6697072652083125

# Row 5 (Sending cluster MuSCs, Target cluster ECs)
$ws.Range("M5").Value = 2.750415333333333
$ws.Range("N5").Value = 8.251245999999998
$ws.Range("O5").Value = 0.04811444325525444
$ws.Range("P5").Value = 0.04811444325525444
$ws.Range("Q5").Value = 0.02173653237933333
$ws.Range("R5").Value = 0.195628791414
$ws.Range("S5").Value = 0.01375220416080564
$ws.Range("T5").Value = 0.01375220416080564

# Row 6 (Sending cluster MuSCs, Target cluster FAPs)
$ws.Range("O6").Value = 0.01415294505639593
$ws.Range("P6").Value = 0.01415294505639593
$ws.Range("S6").Value = 0.004045234169283798
$ws.Range("T6").Value = 0.0040452341692838

# Row 7 (Sending cluster MuSCs, Target cluster MuSCs)
$ws.Range("O7").Value = 0.9377326116883496
$ws.Range("P7").Value = 0.9377326116883496
$ws.Range("S7").Value = 0.2680253464800371
$ws.Range("T7").Value = 0.2680253464800372
